$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.339376211166382
$ws.Range("B1").Value = 2.933096408843994
$ws.Range("C1").Value = 2.759621620178223
$ws.Range("D1").Value = 1.426152944564819
$ws.Range("E1").Value = 1.047430992126465
